# Apply the edit described by the commit:
# "Ongoing QLearning testing (excels) + simplified EtatMap criteria choice with Enum"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RichPhiString($cell, $texts, $subs, $italics) {
    $full = ""
    foreach ($t in $texts) { $full = $full + $t }
    $cell.Value = $full
    $pos = 1
    for ($i = 0; $i -lt $texts.Length; $i++) {
        $len = $texts[$i].Length
        if ($len -gt 0 -and ($subs[$i] -or $italics[$i])) {
            $run = $cell.Characters($pos, $len)
            if ($subs[$i]) { $run.Font.Subscript = $true }
            if ($italics[$i]) { $run.Font.Italic = $true }
        }
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------
# D13: new content -> phi0(s,a), phi1(s,a), phi2(s,a), phi5(s,a)
#      (biais + DistClose + nbGhostDist3)
# ---------------------------------------------------------------------
$texts32 = @(
    "φ", "0", "(s, a), φ", "1", "(s, a), φ", "2", "(s, a), φ5", "(",
    "s, a)`n(", "biais + ", "DistClose", " + nbGhostDist3", ")"
)
$subs32 = @($false,$true,$false,$true,$false,$true,$false,$true,$false,$false,$false,$false,$false)
$italics32 = @($false,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false,$true,$false)

$cellD13 = $ws.Range("D13")
Set-RichPhiString $cellD13 $texts32 $subs32 $italics32

$ws.Range("E13:G13").NumberFormat = "0%"
$ws.Range("E13").Value = 0.85
$ws.Range("F13").Value = 0.69
$ws.Range("G13").Value = 0.01

# I3: "Notes:"
$ws.Range("I3").Value = "Notes:"

# ---------------------------------------------------------------------
# D14: new content -> phi0(s,a), phi1(s,a), phi2(s,a), phi5(s,a)
#      (distanceGhost + nbGhostDist3 +  hasDot + distanceDot)
# ---------------------------------------------------------------------
$texts34 = @(
    "φ", "0", "(s, a), φ", "1", "(s, a), φ", "2", "(s, a), φ5", "(",
    "s, a)`n(distanceGhost", " + nbGhostDist3 +  hasDot + distanceDot", ")"
)
$subs34 = @($false,$true,$false,$true,$false,$true,$false,$true,$false,$false,$false)
$italics34 = @($false,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false)

$ws.Rows(14).RowHeight = 48
$cellD14 = $ws.Range("D14")
Set-RichPhiString $cellD14 $texts34 $subs34 $italics34

$ws.Range("E14:G14").NumberFormat = "0%"
$ws.Range("E14").Value = 0.99
$ws.Range("F14").Value = 0.58
$ws.Range("G14").Value = 0.7

# I5 / I6 notes
$ws.Range("I5").Value = "Manhattan est désatreux"
$ws.Range("I6").Value = "Diviser les distances"

# ---------------------------------------------------------------------
# D15: content that used to be in D13 ("Tous!" variant), moved down one row
# ---------------------------------------------------------------------
$texts29 = @(
    "φ", "0", "(s, a), φ", "1", "(s, a), φ", "2", "(s, a), φ", "3(",
    "s, a),φ4(s, a), φ5(s, a) `n(", "Tous!", ")"
)
$subs29 = @($false,$true,$false,$true,$false,$true,$false,$true,$false,$false,$false)
$italics29 = @($false,$false,$false,$false,$false,$false,$false,$false,$false,$true,$false)

$ws.Rows(15).RowHeight = 33
$cellD15 = $ws.Range("D15")
Set-RichPhiString $cellD15 $texts29 $subs29 $italics29

$ws.Range("E15:G15").NumberFormat = "0%"
$ws.Range("E15").Value = 0.93
$ws.Range("F15").Value = 0.74
$ws.Range("G15").Value = 0.8

# ---------------------------------------------------------------------
# D18 / D19: add empty, right/top aligned + wrap formatted cells
# ---------------------------------------------------------------------
$ws.Range("D18").HorizontalAlignment = -4152
$ws.Range("D18").VerticalAlignment = -4160
$ws.Range("D18").WrapText = $true

$ws.Range("D19").HorizontalAlignment = -4152
$ws.Range("D19").VerticalAlignment = -4160
$ws.Range("D19").WrapText = $true

# ---------------------------------------------------------------------
# Move the " dist(s,a) / (mapX + mapY)" note from D22 to D24
# ---------------------------------------------------------------------
$ws.Range("D22").Clear()
$ws.Range("D24").Value = " dist(s,a) / (mapX + mapY)"

# Remove the now unused formatted-but-empty D25 / D26 cells
$ws.Range("D25").Clear()
$ws.Range("D26").Clear()

# Add new trailing rows 29 / 30 (same style as existing row 28 D28 cell)
$ws.Range("D29").WrapText = $true
$ws.Range("D30").WrapText = $true

# ---------------------------------------------------------------------
# Sheet view: scroll + selection
# ---------------------------------------------------------------------
$ws.Range("J14").Select()
try {
    $excel.ActiveWindow.ScrollRow = 16
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}

Write-Host "Edit applied"
